# Update the financial figures for rows 2-6 (columns D..AJ) and clear the
# now-removed trailing rows 7-9 (columns D..AJ), per the "error solve ifrs
# list" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")

$rowData = @{
  2 = @{ D=18897; E=1081; F=1081; G=-931; H=-1011; I=-1172; J=161; K=40180; L=29238; M=10942; N=4794; O=6149; P=1184; Q=2198; R=-559; S=-1682; T=1002; U=1196; V=19978; W=5.72; X=-5.35; Y=-22.06; Z=-2.44; AA=267.19; AB=969.79; AC=-4949; AD=-2.47; AE=21970; AF=0.56; AG=450; AH=3.69; AI=-8.4; AJ=23206765 }
  3 = @{ D=19271; E=1459; F=1459; G=677; H=447; I=184; J=264; K=40688; L=29573; M=11115; N=4953; O=6161; P=1184; Q=2830; R=-791; S=-1005; T=1208; U=1622; V=19404; W=7.57; X=2.32; Y=3.77; Z=1.11; AA=266.08; AB=981.42; AC=775; AD=17.54; AE=22701; AF=0.6; AG=200; AH=1.47; AI=23.89; AJ=23206765 }
  4 = @{ D=18954; E=1329; F=1329; G=719; H=463; I=279; J=184; K=39671; L=27636; M=12035; N=5313; O=6722; P=1184; Q=1666; R=-89; S=-1485; T=1255; U=411; V=17465; W=7.01; X=2.44; Y=5.44; Z=1.15; AA=229.63; AB=996.09; AC=1179; AD=10.56; AE=24352; AF=0.51; AG=200; AH=1.61; AI=15.71; AJ=23206765 }
  5 = @{ D=18812; E=1001; F=1001; G=195; H=-76; I=-144; J=68; K=40593; L=29244; M=11350; N=4973; O=6377; P=1184; Q=3318; R=-1036; S=-1036; T=931; U=2387; V=16788; W=5.32; X=-0.41; Y=-2.8; Z=-0.19; AA=257.66; AB=973.94; AC=-608; AD=-16.62; AE=22790; AF=0.44; AG=150; AH=1.49; AI=-22.9; AJ=23206765 }
  6 = @{ D=18766; E=1035; F=1035; G=293; H=86; I=-27; K=39893; L=28935; M=10958; N=4842; P=1184; Q=1245; R=-1020; S=-198; T=1205; U=41; V=16902; W=5.52; X=0.46; Y=-0.55; Z=0.22; AA=264.06; AB=959.43; AC=-115; AD=-61.14; AE=22176; AF=0.32; AG=200; AH=2.85; AI=-161.72; AJ=23206765 }
}

foreach ($r in $rowData.Keys) {
  $cols = $rowData[$r]
  foreach ($col in $headers) {
    if ($cols.ContainsKey($col)) {
      $ws.Range("$col$r").Value = $cols[$col]
    }
  }
}

# Rows 7-9 lost all of their financial-figure columns (D..AJ); only the
# leading identifier columns A-C remain populated.
$ws.Range("D7:AJ9").ClearContents()
